$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo" shifts from D to E)
$ws.Columns.Item(4).Insert()

# Apply the same header formatting used by the other header cells
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Update existing values per diff
$ws.Range("B2").Value = 0.4936382884411732
$ws.Range("C2").Value = 0.9901706133763721
$ws.Range("D2").Value = 0.575257090994669

$ws.Range("A1:E2").Select()
